# Apply the crypto price / volume updates for Wed Jun 12 13:25:49 UTC 2024 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.320.50'
$ws.Range("E2").Value = '  +3.07%  '
$ws.Range("D3").Value = '3.632.36'
$ws.Range("E3").Value = '  +2.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '627.33'
$ws.Range("E5").Value = '  +2.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.01'
$ws.Range("E6").Value = '  +3.76%  '
$ws.Range("D7").Value = '3.630.80'
$ws.Range("E7").Value = '  +2.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  +1.79%  '
$ws.Range("E10").Value = '  +2.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.19'
$ws.Range("E11").Value = '  +4.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("E12").Value = '  +2.32%  '
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.27'
$ws.Range("E14").Value = '  +4.03%  '
$ws.Range("D15").Value = '4.251.09'
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").Value = '3.628.79'
$ws.Range("E16").Value = '  +2.49%  '
$ws.Range("D17").Value = '69.351.82'
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("E19").Value = '  +4.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.87'
$ws.Range("E20").Value = '  +2.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.24'
$ws.Range("E21").Value = '  +10.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.62'
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.642'
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.57'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000137'
$ws.Range("E25").Value = '  +10.97%  '
$ws.Range("D26").Value = '3.781.60'
$ws.Range("E26").Value = '  +2.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.64'
$ws.Range("E27").Value = '  +3.78%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.12'
$ws.Range("E29").Value = '  +10.23%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.62'
$ws.Range("E30").Value = '  +3.06%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.71'
$ws.Range("E31").Value = '  +2.01%  '
$ws.Range("E32").Value = '  +10.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.62'
$ws.Range("E33").Value = '  +6.99%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  +4.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.41'
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("D37").Value = '3.622.95'
$ws.Range("E37").Value = '  +2.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.38'
$ws.Range("E38").Value = '  +4.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.35'
$ws.Range("E39").Value = '  +9.14%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0927'
$ws.Range("E41").Value = '  +6.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '175.42'
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.61'
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.67'
$ws.Range("E45").Value = '  +14.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.916'
$ws.Range("E46").Value = '  +2.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.37'
$ws.Range("E47").Value = '  +11.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("E48").Value = '  +7.47%  '
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.79'
$ws.Range("E50").Value = '  +2.68%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.268'
$ws.Range("E51").Value = '  +7.35%  '
